# Adds the 2019 annual data column (AD) to the OOSS strikes table,
# as described in the commit "se agregan datos del anuario 2019 para OOSS y huelgas".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AD (year 2019), following the existing year headers in row 1.
# Leading apostrophe forces the value to be stored as text (matching the other
# year headers in row 1, e.g. B1:AC1 = "1991".."2018") instead of being
# auto-coerced to a number.
$ws.Range("AD1").Value = "'2019"

# New 2019 values per activity/branch row (column AD), in row order 2..12.
$ws.Range("AD2").Value  = 1562   # Actividades no especiftcadas
$ws.Range("AD3").Value  = 1260   # Agricultura y pesca
$ws.Range("AD4").Value  = 1960   # Comercio
$ws.Range("AD5").Value  = 335    # Construcción
$ws.Range("AD6").Value  = 214    # Electricidad, gas y agua
$ws.Range("AD7").Value  = 183    # Establecimientos ftnancieros
$ws.Range("AD8").Value  = 1365   # Industria
$ws.Range("AD9").Value  = 298    # Minería
$ws.Range("AD10").Value = 3004   # Servicios
$ws.Range("AD11").Value = 11926  # Total
$ws.Range("AD12").Value = 1745   # Transporte y comunicaciones
